# feat: add 2022-Q1 data
#
# The workbook has two sheets: "2021-Q1" (per-quarter holdings) and "总计"
# (summary roll-up). This adds a new "2022-Q1" holdings sheet (positioned
# between the two existing sheets) and records its summary row at the top
# of "总计" (pushing the existing "2021-Q1" summary row down).

$wb = $excel.ActiveWorkbook
$wsQ1 = $wb.Worksheets.Item(1)          # "2021-Q1"
$wsTotal = $wb.Worksheets.Item("总计")  # "总计"

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet by duplicating "总计" (so it starts
#    with the same header/body cell styling used on that sheet), then
#    move it right after "2021-Q1" and rename it.
# ---------------------------------------------------------------------
$wsTotal.Copy($null, $wsQ1)
$wsNew = $wb.Worksheets.Item(2)
$wsNew.Name = "2022-Q1"

# NOTE: after Worksheet.Copy(), the original $wsTotal handle resolves to
# the newly-created copy rather than the source sheet, so re-fetch a
# fresh handle to the real "总计" sheet for use below.
$wsTotal = $wb.Worksheets.Item("总计")

# Stretch the styled header (B1:D1) across to H1, and the styled body
# row (B2:D2) across to H2, by copying the formatting of column D.
$wsNew.Range("D1").Copy($wsNew.Range("E1"))
$wsNew.Range("D1").Copy($wsNew.Range("F1"))
$wsNew.Range("D1").Copy($wsNew.Range("G1"))
$wsNew.Range("D1").Copy($wsNew.Range("H1"))

$wsNew.Range("D2").Copy($wsNew.Range("E2"))
$wsNew.Range("D2").Copy($wsNew.Range("F2"))
$wsNew.Range("D2").Copy($wsNew.Range("G2"))
$wsNew.Range("D2").Copy($wsNew.Range("H2"))

# Replicate the A2:H2 row layout down through row 8 (7 fund rows total).
for ($r = 3; $r -le 8; $r++) {
    $wsNew.Range("A2:H2").Copy($wsNew.Range("A" + $r + ":H" + $r))
}

# Header row text.
$wsNew.Range("B1").Value = "基金代码"
$wsNew.Range("C1").Value = "基金名称"
$wsNew.Range("D1").Value = "基金规模"
$wsNew.Range("E1").Value = "股票总仓位"
$wsNew.Range("F1").Value = "仓位占比"
$wsNew.Range("G1").Value = "持有市值(亿元)"
$wsNew.Range("H1").Value = "仓位排名"

# Fund holding rows. Columns B-G keep their original text formatting
# (fund codes such as "012073" must not collapse to numbers), column H
# (仓位排名) is numeric.
$rows = @(
    @("012073", "华安均衡优选混合A", "8.33", "89.23", "2.05", "0.1708", 10),
    @("001581", "华安沪港深通精选灵活配置混合", "4.92", "92.91", "3.13", "0.1540", 9),
    @("040018", "华安香港精选股票(QDII)", "5.47", "88.46", "2.41", "0.1318", 10),
    @("011144", "华安汇宏精选混合A", "1.07", "85.87", "4.20", "0.0449", 2),
    @("011145", "华安汇宏精选混合C", "0.25", "85.87", "4.20", "0.0105", 2),
    @("040021", "华安大中华升级股票(QDII)", "0.26", "87.37", "2.31", "0.0060", 9),
    @("012074", "华安均衡优选混合C", "0.25", "89.23", "2.05", "0.0051", 10)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = 2 + $i
    $data = $rows[$i]

    $wsNew.Range("A" + $r).Value = $i

    $codeCell = $wsNew.Range("B" + $r)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $data[0]

    $wsNew.Range("C" + $r).Value = $data[1]

    $col3 = $wsNew.Range("D" + $r)
    $col3.NumberFormat = "@"
    $col3.Value = $data[2]

    $col4 = $wsNew.Range("E" + $r)
    $col4.NumberFormat = "@"
    $col4.Value = $data[3]

    $col5 = $wsNew.Range("F" + $r)
    $col5.NumberFormat = "@"
    $col5.Value = $data[4]

    $col6 = $wsNew.Range("G" + $r)
    $col6.NumberFormat = "@"
    $col6.Value = $data[5]

    $wsNew.Range("H" + $r).Value = $data[6]
}

# ---------------------------------------------------------------------
# 2) Insert the 2022-Q1 summary row at the top of "总计", pushing the
#    existing 2021-Q1 summary row down to row 3.
# ---------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("A2:D2").ClearFormats()
$wsTotal.Range("A3").Copy($wsTotal.Range("A2"))

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 7
$wsTotal.Range("D2").Value = 0.52

$wsTotal.Range("A3").Value = 1
